$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing row 16 values (revised measurements) ---
$ws.Range("D16").Value = 10.199999999999999
$ws.Range("E16").Value = 1013.3
$ws.Range("F16").Value = 76

# --- Append the new row 17: May 2022 data ---
$ws.Range("A17").Value = 2022
$ws.Range("B17").Value = 5
$ws.Range("C17").Formula = "=_xlfn.CONCAT(A17,`" `",B17)"
$ws.Range("D17").Value = 18.5
$ws.Range("E17").Value = 1017.9
$ws.Range("F17").Value = 23

# --- Extend each chart series so it covers the new row ---
# Chart 1 plots column D (temperature), Chart 2 plots column E (pressure),
# Chart 3 plots column F (humidity); all share column C as the category axis.
$charts = $ws.ChartObjects()
$charts.Item(1).Chart.SeriesCollection(1).Formula = "=SERIES(,monthlyData!`$C`$1:`$C`$17,monthlyData!`$D`$1:`$D`$17,1)"
$charts.Item(2).Chart.SeriesCollection(1).Formula = "=SERIES(,monthlyData!`$C`$1:`$C`$17,monthlyData!`$E`$1:`$E`$17,1)"
$charts.Item(3).Chart.SeriesCollection(1).Formula = "=SERIES(,monthlyData!`$C`$1:`$C`$17,monthlyData!`$F`$1:`$F`$17,1)"

# --- Remove the stale hidden chart-helper defined names (no longer match the range) ---
$wb.Names.Item("_xlchart.v1.0").Delete()
$wb.Names.Item("_xlchart.v1.1").Delete()
$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

# --- Put the selection where the user left off after typing the new row ---
$ws.Range("F17").Select() | Out-Null
